# Apply updated voltage magnitude (vm_pu) results for the 380 kV case.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.020715172285117
$ws.Cells.Item(2, 4).Value = 1.029780228227659
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.036481785142401
$ws.Cells.Item(2, 9).Value = 1.029393502144703
$ws.Cells.Item(2, 10).Value = 1.025910718035589
$ws.Cells.Item(2, 11).Value = 1.032593146699511
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.039275401874681
$ws.Cells.Item(2, 14).Value = 1.012715097980796
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.021518409425796
$ws.Cells.Item(3, 4).Value = 1.030380914748545
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.037269254840939
$ws.Cells.Item(3, 9).Value = 1.029479104628175
$ws.Cells.Item(3, 10).Value = 1.026352024954584
$ws.Cells.Item(3, 11).Value = 1.033002608787188
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.039872540812007
$ws.Cells.Item(3, 14).Value = 1.012861942845745
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.02203877066807
$ws.Cells.Item(4, 4).Value = 1.030769959850846
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.03777967932735
$ws.Cells.Item(4, 9).Value = 1.0295332420575
$ws.Cells.Item(4, 10).Value = 1.026637540964479
$ws.Cells.Item(4, 11).Value = 1.033267223135638
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.040259151714285
$ws.Cells.Item(4, 14).Value = 1.012956918131548
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.022257675571742
$ws.Cells.Item(5, 4).Value = 1.030933598857961
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.037994470096103
$ws.Cells.Item(5, 9).Value = 1.029555700747461
$ws.Cells.Item(5, 10).Value = 1.026757561239892
$ws.Cells.Item(5, 11).Value = 1.033378385397541
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.040421734348042
$ws.Cells.Item(5, 14).Value = 1.012996834944017
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.022294439122798
$ws.Cells.Item(6, 4).Value = 1.030961079463365
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.038030546565751
$ws.Cells.Item(6, 9).Value = 1.02955945400793
$ws.Cells.Item(6, 10).Value = 1.026777712509597
$ws.Cells.Item(6, 11).Value = 1.033397045203757
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.040449035651285
$ws.Cells.Item(6, 14).Value = 1.013003536505763
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.022041695118387
$ws.Cells.Item(7, 4).Value = 1.030772146073874
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.037782548556229
$ws.Cells.Item(7, 9).Value = 1.029533543333965
$ws.Cells.Item(7, 10).Value = 1.026639144724775
$ws.Cells.Item(7, 11).Value = 1.033268708813669
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.040261323951596
$ws.Cells.Item(7, 14).Value = 1.012957451544962
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.020986502044061
$ws.Cells.Item(8, 4).Value = 1.029983157199668
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.036747730960146
$ws.Cells.Item(8, 9).Value = 1.029422690952672
$ws.Cells.Item(8, 10).Value = 1.026059867084568
$ws.Cells.Item(8, 11).Value = 1.032731594846348
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.039477160735849
$ws.Cells.Item(8, 14).Value = 1.012764733547933
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.019131889654449
$ws.Cells.Item(9, 4).Value = 1.028595705735653
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.034931074728898
$ws.Cells.Item(9, 9).Value = 1.029217788154204
$ws.Cells.Item(9, 10).Value = 1.025038864048028
$ws.Cells.Item(9, 11).Value = 1.03178262797819
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.038097143425646
$ws.Cells.Item(9, 14).Value = 1.012424829780809
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.017898793968161
$ws.Cells.Item(10, 4).Value = 1.027672767264229
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.033724684625465
$ws.Cells.Item(10, 9).Value = 1.029074795343038
$ws.Cells.Item(10, 10).Value = 1.024358110269604
$ws.Cells.Item(10, 11).Value = 1.031148384298291
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.037178432379706
$ws.Cells.Item(10, 14).Value = 1.012198046512325
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.017365657470499
$ws.Cells.Item(11, 4).Value = 1.027273630241966
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.033203449065946
$ws.Cells.Item(11, 9).Value = 1.029011370713859
$ws.Cells.Item(11, 10).Value = 1.024063332545104
$ws.Cells.Item(11, 11).Value = 1.030873388255704
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.036780950956363
$ws.Cells.Item(11, 14).Value = 1.012099809825894
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.017167749097803
$ws.Cells.Item(12, 4).Value = 1.027125450280844
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.033010012297803
$ws.Cells.Item(12, 9).Value = 1.028987586137449
$ws.Cells.Item(12, 10).Value = 1.023953839262737
$ws.Cells.Item(12, 11).Value = 1.030771189172851
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.036633359505186
$ws.Cells.Item(12, 14).Value = 1.01206331510336
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.017210195554773
$ws.Cells.Item(13, 4).Value = 1.02715723185597
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.033051497263201
$ws.Cells.Item(13, 9).Value = 1.028992698218805
$ws.Cells.Item(13, 10).Value = 1.023977325909241
$ws.Cells.Item(13, 11).Value = 1.030793113612426
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.036665016042794
$ws.Cells.Item(13, 14).Value = 1.012071143574939
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.017349295798423
$ws.Cells.Item(14, 4).Value = 1.02726138005189
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.033187455974071
$ws.Cells.Item(14, 9).Value = 1.029009409276566
$ws.Cells.Item(14, 10).Value = 1.024054281782704
$ws.Cells.Item(14, 11).Value = 1.030864941526848
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.036768749953348
$ws.Cells.Item(14, 14).Value = 1.012096793264983
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.017435016302741
$ws.Cells.Item(15, 4).Value = 1.027325559494541
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.033271247654989
$ws.Cells.Item(15, 9).Value = 1.029019675604592
$ws.Cells.Item(15, 10).Value = 1.024101696913971
$ws.Cells.Item(15, 11).Value = 1.030909190054397
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.036832670630587
$ws.Cells.Item(15, 14).Value = 1.012112596204868
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.017934193530438
$ws.Cells.Item(16, 4).Value = 1.027699267406345
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.033759301501193
$ws.Cells.Item(16, 9).Value = 1.029078972917638
$ws.Cells.Item(16, 10).Value = 1.024377673664651
$ws.Cells.Item(16, 11).Value = 1.031166627325776
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.037204818945082
$ws.Cells.Item(16, 14).Value = 1.012204565396551
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.018247530195301
$ws.Cells.Item(17, 4).Value = 1.02793382014218
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.034065751335753
$ws.Cells.Item(17, 9).Value = 1.029115765358222
$ws.Cells.Item(17, 10).Value = 1.024550785687111
$ws.Cells.Item(17, 11).Value = 1.031328014389753
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.037438346435063
$ws.Cells.Item(17, 14).Value = 1.012262245390538
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.018430371440963
$ws.Cells.Item(18, 4).Value = 1.028070679033384
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.034244608122493
$ws.Cells.Item(18, 9).Value = 1.029137080239745
$ws.Cells.Item(18, 10).Value = 1.024651758252116
$ws.Cells.Item(18, 11).Value = 1.031422113484027
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.037574590453652
$ws.Cells.Item(18, 14).Value = 1.012295885453824
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.018492728633608
$ws.Cells.Item(19, 4).Value = 1.028117352502447
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.034305612220124
$ws.Cells.Item(19, 9).Value = 1.029144323358863
$ws.Cells.Item(19, 10).Value = 1.024686187120277
$ws.Cells.Item(19, 11).Value = 1.031454192820939
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.037621051411991
$ws.Cells.Item(19, 14).Value = 1.012307355203913
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.018213904131434
$ws.Cells.Item(20, 4).Value = 1.027908649848626
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.034032860797785
$ws.Cells.Item(20, 9).Value = 1.029111832925108
$ws.Cells.Item(20, 10).Value = 1.024532212473364
$ws.Cells.Item(20, 11).Value = 1.031310702720859
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.037413287886283
$ws.Cells.Item(20, 14).Value = 1.012256057252508
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.017308330854536
$ws.Cells.Item(21, 4).Value = 1.027230708854371
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.033147414720426
$ws.Cells.Item(21, 9).Value = 1.029004494515502
$ws.Cells.Item(21, 10).Value = 1.024031620177475
$ws.Cells.Item(21, 11).Value = 1.030843791456414
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.036738201484183
$ws.Cells.Item(21, 14).Value = 1.012089240211526
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.016739669226092
$ws.Cells.Item(22, 4).Value = 1.026804908651131
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.032591703508934
$ws.Cells.Item(22, 9).Value = 1.028935700265695
$ws.Cells.Item(22, 10).Value = 1.023716880258963
$ws.Cells.Item(22, 11).Value = 1.030549918092966
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.036314042940867
$ws.Cells.Item(22, 14).Value = 1.011984325613953
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.017041059794902
$ws.Cells.Item(23, 4).Value = 1.027030590254847
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.032886200660626
$ws.Cells.Item(23, 9).Value = 1.028972292995662
$ws.Cells.Item(23, 10).Value = 1.023883729184569
$ws.Cells.Item(23, 11).Value = 1.030705734614518
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.036538868843172
$ws.Cells.Item(23, 14).Value = 1.012039945523499
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.018229098063743
$ws.Cells.Item(24, 4).Value = 1.027920023072016
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.034047722277204
$ws.Cells.Item(24, 9).Value = 1.0291136102721
$ws.Cells.Item(24, 10).Value = 1.024540604911662
$ws.Cells.Item(24, 11).Value = 1.031318525227956
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.037424610669099
$ws.Cells.Item(24, 14).Value = 1.012258853417215
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.019610775149119
$ws.Cells.Item(25, 4).Value = 1.028954045909203
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.035399902211755
$ws.Cells.Item(25, 9).Value = 1.029271889903556
$ws.Cells.Item(25, 10).Value = 1.025302838473076
$ws.Cells.Item(25, 11).Value = 1.032028246660292
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.038453689420085
$ws.Cells.Item(25, 14).Value = 1.012512736936458
